$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the most recent existing data row (139) with the refreshed figures
$ws.Cells.Item(139, 2).Value = 6
$ws.Cells.Item(139, 3).Value = 2.1
$ws.Cells.Item(139, 4).Value = 6.4

# Append the new month (01-07-2021) as row 140.
# Enter the date label via a text formula first so Excel's General-format
# date auto-recognition doesn't kick in, then collapse the formula down to
# its static text result in place (paste values only) - this keeps column A
# as text, matching the rest of the "Serie" column, without leaving behind
# any formula or any new/unused cell style.
$ws.Cells.Item(140, 1).Formula = "=""01-07-2021"""
$ws.Cells.Item(140, 1).Copy()
$ws.Cells.Item(140, 1).PasteSpecial(-4163)

$ws.Cells.Item(140, 2).Value = 6.5
$ws.Cells.Item(140, 3).Value = 1.9
$ws.Cells.Item(140, 4).Value = 6.6
